$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append four new SKU rows (88-91) below the existing data (A1:A87).
$ws.Cells.Item(88, 1).Value = 10032499
$ws.Cells.Item(89, 1).Value = 10007485
$ws.Cells.Item(90, 1).Value = 10105349
$ws.Cells.Item(91, 1).Value = 10220817

# Rows 89 and 91 carry a distinct look: a small grey "Open Sans" font and a
# taller row height, matching the two "highlighted" rows added upstream.
$r89 = $ws.Range("A89")
$r89.Font.Name = "Open Sans"
$r89.Font.Size = 9
$r89.Font.Color = 4473924

# Re-use the exact same formatting for row 91 via Copy/PasteSpecial so both
# rows end up sharing a single cell style instead of each minting its own.
$r89.Copy()
$ws.Range("A91").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(89).RowHeight = 15.75
$ws.Rows.Item(91).RowHeight = 15.75

# Restore the view: scrolled down to row 71, with F77 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 71
$ws.Range("F77").Select()
